# ** Update ** All features are integrated, except for noise floor.
#
# Adds a "percent change" summary block (rows 10-16) under the existing
# Bx/By/Bz/H/D/I/F comparison table: for each quantity, how much the WMM
# (col P) differs from the generated/NOAA value (col N), and how much the
# measured/calculated value (col R) differs from WMM (col P).
# Also updates the sheet's view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- View state: scrolled so column I is the left-most visible column,
# with A5 selected as the active cell.
$win = $excel.ActiveWindow
$win.ScrollColumn = 9
$win.ScrollRow = 1
$ws.Range("A5").Select()

# --- New rows 10-16: row labels in col N (same strings/order as A2:A8,
# just re-ordered) with the centred style used by the existing A-column
# labels, then the two percentage-difference formulas in columns P and R.
$labels = @("D", "I", "H", "By", "Bx", "Bz", "F")
for ($i = 0; $i -lt $labels.Length; $i++) {
    $destRow = 10 + $i
    $srcRow  = 2 + $i

    $ws.Range("N$destRow").Value = $labels[$i]

    $ws.Range("P$destRow").Formula = "=(P$srcRow-N$srcRow)/N$srcRow * 100"
    $ws.Range("R$destRow").Formula = "=(R$srcRow-P$srcRow)/P$srcRow * 100"

    # Formula entry auto-inherits the custom "0.000" number format from the
    # referenced N/P/R precedent cells - strip that back off so P/R keep the
    # workbook's default (General) style, matching the rest of the summary.
    $ws.Range("P$destRow").ClearFormats()
    $ws.Range("R$destRow").ClearFormats()
}

# Centre-align the new N10:N16 labels (same look as A2:A8).
$ws.Range("N10:N16").HorizontalAlignment = -4108  # xlCenter
